# Rename model Region to Zone for clarity
$wb = $excel.ActiveWorkbook

# Rename the "Region" worksheet to "Zone"
$ws = $wb.Worksheets.Item("Region")
$ws.Name = "Zone"

# Select/activate the renamed sheet so it becomes the active tab
$ws.Activate()
$ws.Select()
